$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New client / service row appended to the revenue report (row 30)
$ws.Range("A30").Value = "Dawson"
$ws.Range("B30").Value = "Suivi - Accompagnement spirituel et énergétique :: Follow-up session - Spiritual and Energetic Guidance"
$ws.Range("C30").Value = 100
$ws.Range("D30").Value = 1
$ws.Range("F30").Value = 100
$ws.Range("H30").Value = 100

# Reflect the user's final selection after the edit (sheet stays frozen at row 1)
$null = $ws.Range("B36").Select()
